$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.00340530918659546
$ws.Range("C2").Value = 0.000735237210742203
$ws.Range("D2").Value = 0.00065784382013776
$ws.Range("E2").Value = 0.00537884064700875
$ws.Range("F2").Value = 0.0219797229316616
$ws.Range("G2").Value = 0.00150917111678663
$ws.Range("H2").Value = 0.00247658849934216
$ws.Range("I2").Value = 0.027822923922297
$ws.Range("J2").Value = 0.00944199365374197
$ws.Range("K2").Value = 0.00232180171813327
$ws.Range("L2").Value = 0.000386966953022212
$ws.Range("M2").Value = 0.00371488274901323
$ws.Range("N2").Value = 0.004024456311431
$ws.Range("O2").Value = 0.906315300673322
$ws.Range("P2").Value = 0.00100611407785775
$ws.Range("Q2").Value = 0.000580450429533318
$ws.Range("R2").Value = 0.000580450429533318
$ws.Range("S2").Value = 0.985837009519387
$ws.Range("T2").Value = 0.000464360343626654
$ws.Range("U2").Value = 0.0012769909449733
$ws.Range("V2").Value = 0.922103552356629
$ws.Range("W2").Value = 0.00417924309263989
$ws.Range("X2").Value = 0.0272811701880659
$ws.Range("B3").Value = 0.985566132652271
$ws.Range("C3").Value = 0.991061063385187
$ws.Range("D3").Value = 0.00448881665505766
$ws.Range("E3").Value = 0.00243789180403994
$ws.Range("F3").Value = 0.0139308103087996
$ws.Range("G3").Value = 0.00139308103087996
$ws.Range("H3").Value = 0.965598637876325
$ws.Range("I3").Value = 0.969932667750174
$ws.Range("J3").Value = 0.000619147124835539
$ws.Range("K3").Value = 0.905386579986069
$ws.Range("L3").Value = 0.00147047442148441
$ws.Range("M3").Value = 0.000580450429533318
$ws.Range("N3").Value = 0.993808528751645
$ws.Range("O3").Value = 0.0769677269561179
$ws.Range("P3").Value = 0.00506926708459098
$ws.Range("Q3").Value = 0.955614890488352
$ws.Range("R3").Value = 0.000619147124835539
$ws.Range("S3").Value = 0.000696540515439981
$ws.Range("T3").Value = 0.994737249438898
$ws.Range("U3").Value = 0.991022366689885
$ws.Range("V3").Value = 0.00352139927250213
$ws.Range("W3").Value = 0.959561953409179
$ws.Range("X3").Value = 0.968230013156876
$ws.Range("B4").Value = 0.00096741738255553
$ws.Range("C4").Value = 0.000270876867115548
$ws.Range("D4").Value = 0.00123829424967108
$ws.Range("E4").Value = 0.986262673167711
$ws.Range("F4").Value = 0.00472099682687099
$ws.Range("G4").Value = 0.816616360962774
$ws.Range("H4").Value = 0.00464360343626654
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.00236049841343549
$ws.Range("K4").Value = 0.00576580760003096
$ws.Range("L4").Value = 0.0000386966953022212
$ws.Range("M4").Value = 0.00116090085906664
$ws.Range("N4").Value = 0.00065784382013776
$ws.Range("O4").Value = 0.00499187369398653
$ws.Range("P4").Value = 0.862162371333488
$ws.Range("Q4").Value = 0.00344400588189769
$ws.Range("R4").Value = 0.0071201919356087
$ws.Range("S4").Value = 0.0100611407785775
$ws.Range("T4").Value = 0.00394706292082656
$ws.Range("U4").Value = 0.00325052240538658
$ws.Range("V4").Value = 0.0707762557077626
$ws.Range("W4").Value = 0.0112607383329464
$ws.Range("X4").Value = 0.00297964553827103
$ws.Range("B5").Value = 0.0098676573020664
$ws.Range("C5").Value = 0.00770064236514202
$ws.Range("D5").Value = 0.993537651884529
$ws.Range("E5").Value = 0.00572711090472874
$ws.Range("F5").Value = 0.959020199674948
$ws.Range("G5").Value = 0.178662642210355
$ws.Range("H5").Value = 0.026584629672626
$ws.Range("I5").Value = 0.00170265459329773
$ws.Range("J5").Value = 0.987307483940871
$ws.Range("K5").Value = 0.0847070660165622
$ws.Range("L5").Value = 0.998065165234889
$ws.Range("M5").Value = 0.994272889095271
$ws.Range("N5").Value = 0.00147047442148441
$ws.Range("O5").Value = 0.0112994350282486
$ws.Range("P5").Value = 0.130640043340299
$ws.Range("Q5").Value = 0.0396641126847767
$ws.Range("R5").Value = 0.991680210510022
$ws.Range("S5").Value = 0.00282485875706215
$ws.Range("T5").Value = 0.000851327296648866
$ws.Range("U5").Value = 0.00445011995975544
$ws.Range("V5").Value = 0.0027474653664577
$ws.Range("W5").Value = 0.0243402213450971
$ws.Range("X5").Value = 0.00123829424967108
Write-Output "done"
